$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Sheet is protected - unprotect to allow edits, then restore protection after
$ws.Unprotect()

# Update the confidential disclosure date (A13): 2021-06-09 -> 2021-06-10
$ws.Range("A13").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-06-10 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) columns for rows 2-10
$ws.Range("D2").Value = 0.1319144550555642
$ws.Range("E2").Value = 0.02705564051288056

$ws.Range("D3").Value = 0.1089618002823754
$ws.Range("E3").Value = -0.001644350481559753

$ws.Range("D4").Value = 0.1101291338750764
$ws.Range("E4").Value = 0.002848564169682044

$ws.Range("D5").Value = 0.1190526689078917
$ws.Range("E5").Value = 0.001246979970384254

$ws.Range("D6").Value = 0.1201529732042597
$ws.Range("E6").Value = 0.005485463521667455

$ws.Range("D7").Value = 0.149894463892681
$ws.Range("E7").Value = -0.005571030640668551

$ws.Range("D8").Value = 0.130152336441942
$ws.Range("E8").Value = 0.006530380465644692

$ws.Range("D9").Value = 0.1297421683402096
$ws.Range("E9").Value = 0.01299141262011871

$ws.Range("E10").Value = 0.00621153130292762

# Restore sheet protection
$ws.Protect()
